$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Simple, same-run text replacements -------------------------------------------------

Replace-Text "Unveiling the Mysteries of Dark Matter" "Biology: The Tapestry of Life"
Replace-Text " Alana Hendricks" " Sarah Johnson"
Replace-Text "alanahendricks@astronomysociety" "sarahj@educatormail"

Replace-Text "For decades, astronomers and physicists have been engrossed in a captivating cosmic enigma: the existence and nature of dark matter" "Biology, the study of life, is an awe-inspiring exploration into the diverse tapestry of organisms that inhabit our planet"
Replace-Text " This enigmatic substance, believed to comprise over 26% of the universe, exerts a gravitational influence far exceeding its apparent mass, shaping galaxies, and influencing the universe's expansion rate" " It captivates us with its complexities and unravels the intricate web that connects all living things"
Replace-Text " Yet, despite its profound impact, dark matter remains shrouded in mystery, fueling scientific curiosity and propelling groundbreaking research" " Embark on a journey to discover the vibrant ecosystems that thrive in harmony, the astonishing adaptations that organisms possess, and the profound impact we have on the delicate balance of life"

Replace-Text "The evidence for dark matter's existence is compelling" "As we delve into the microscopic realm of cells, we uncover astonishing marvels of organization and functionality"
Replace-Text " Through meticulous observations of galaxies and galaxy clusters, astronomers have discovered that the gravitational force necessary to hold these celestial structures together far surpasses the gravitational pull exerted by the visible matter they contain" " Each cell is an intricate microcosm within itself, carrying out essential processes that sustain life"
Replace-Text " This discrepancy suggests the presence of an invisible mass, an unseen entity governing the universe's dynamics" " The diversity of life becomes evident as we encounter the remarkable variations in form and behavior among organisms"

Replace-Text "Furthermore, observations of the cosmic microwave background radiation, the leftover glow from the Big Bang, provide further clues about dark matter's existence" "Biology enables us to comprehend how organisms interact with each other and with their surroundings, forming complex ecosystems that thrive through interconnectedness"
Replace-Text " Minute temperature variations in this radiation hint at the gravitational influence of dark matter during the universe's early moments, supporting the notion that it played a pivotal role in shaping the universe's structure" " We learn how delicate balances are maintained within these ecosystems and the intricate roles that each organism plays in preserving this equilibrium"

Replace-Text "The quest to understand dark matter has captivated scientists worldwide, driving cutting-edge research and groundbreaking discoveries" "Biology is an intriguing subject that unveils the mysteries of life's tapestry"
Replace-Text " While its true identity remains elusive, the evidence for its existence is undeniable" " It encompasses the study of cells, their intricate organization and functionality, the bewildering diversity of organisms, and the dynamic interactions between organisms within ecosystems"
Replace-Text " Through continued exploration and innovation, scientists are determined to unravel the mysteries of dark matter, shedding light on one of the universe's most enigmatic components" " Biology nurtures an understanding of the influence we have on the environment and inspires us to act as responsible stewards of our planet"

# --- Splits: a trailing "." run gets duplicated into its own run, followed by a brand-new
#     sentence run. We locate the sentence that now immediately precedes the lone "." run,
#     then replace that single "." character with "." + new sentence as two distinct runs
#     (matching the rPr of the surrounding text exactly). --------------------------------

function Split-AfterSentence($precedingText, $newSentence, $sz) {
    $rng = $d.Content
    $rng.Find.Execute($precedingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $periodStart = $rng.End
    $periodRange = $d.Range($periodStart, $periodStart + 1)

    if ($sz) {
        $szXml = "<w:sz w:val=""$sz""/>"
    } else {
        $szXml = ""
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/>' + $szXml + '</w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/>' + $szXml + '</w:rPr><w:t xml:space="preserve"> ' + $newSentence + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $periodRange.InsertXML($xml) | Out-Null
}

Split-AfterSentence " The diversity of life becomes evident as we encounter the remarkable variations in form and behavior among organisms" "From the grandeur of the majestic whales that roam our oceans to the minuscule yet tenacious microorganisms, we marvel at the myriad life forms that grace our planet" "24"
Split-AfterSentence " We learn how delicate balances are maintained within these ecosystems and the intricate roles that each organism plays in preserving this equilibrium" "The interdependence of organisms within these interconnected systems highlights the profound responsibility we bear as stewards of our natural world" "24"
Split-AfterSentence " Biology nurtures an understanding of the influence we have on the environment and inspires us to act as responsible stewards of our planet" "As we continue to explore the wonders of life, we unravel the secrets of our existence and uncover the boundless possibilities that the study of biology holds" $null

# --- Append a trailing empty paragraph --------------------------------------------------

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null
